# Add team record (Wins/Losses/Ties) columns AD:AF to the LAD_2021 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - matches style of existing header cells (A1:AC1):
# bold font, centered horizontal/top vertical alignment, thin border all around.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRng = $ws.Range("AD1:AF1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108   # xlCenter
$headerRng.VerticalAlignment = -4160     # xlTop
$headerRng.Borders.LineStyle = 1         # xlContinuous (thin, all sides)

# Every player row (2-63) gets the same team record: 106 wins, 56 losses, 0 ties
$lastRow = 63
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 106   # AD
    $ws.Cells.Item($r, 31).Value = 56    # AE
    $ws.Cells.Item($r, 32).Value = 0     # AF
}
